$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2 (G2 = 5489)
$ws.Range("H2").Value = 101.25
$ws.Range("I2").Value = 103.333336
$ws.Range("J2").Value = 95
$ws.Range("K2").Value = 103.333336
$ws.Range("L2").Value = 95
$ws.Range("M2").Value = 9.666663999999997
$ws.Range("N2").Value = -321

# Row 33 (G33 = 5512)
$ws.Range("H33").Value = 134.06667
$ws.Range("I33").Value = 144.3077
$ws.Range("J33").Value = 67.5
$ws.Range("K33").Value = 144.3077
$ws.Range("L33").Value = 67.5
$ws.Range("M33").Value = 84.69229999999999
$ws.Range("N33").Value = -525.5

# Row 92 (G92 = 19901)
$ws.Range("H92").Value = 587.125
$ws.Range("I92").Value = 339.8
$ws.Range("K92").Value = 339.8
$ws.Range("M92").Value = 908.2

# Row 98 (G98 = 36237)
$ws.Range("H98").Value = 1213.3334
$ws.Range("I98").Value = 820
$ws.Range("K98").Value = 820
$ws.Range("M98").Value = 678

# Row 122 (G122 = 36237)
$ws.Range("H122").Value = 1213.3334
$ws.Range("I122").Value = 820
$ws.Range("K122").Value = 2460
$ws.Range("M122").Value = -10

# Row 134 (G134 = 41997)
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

# Row 135 (G135 = 44047)
$ws.Range("H135").Value = 899
$ws.Range("I135").Value = 848.9167
$ws.Range("K135").Value = 7640.2503
$ws.Range("M135").Value = -5105.2503

# Row 137 (G137 = 44013)
$ws.Range("H137").Value = 2590
$ws.Range("I137").Value = 2350.75
$ws.Range("J137").Value = 3000.1428
$ws.Range("K137").Value = 7052.25
$ws.Range("L137").Value = 9000.428400000001
$ws.Range("M137").Value = -4502.25
$ws.Range("N137").Value = -14100.4284

$ws = $wb.Worksheets.Item("ARM")
# Row 61 (G61 = 43999)
$ws.Range("H61").Value = 1598.2
$ws.Range("I61").Value = 1498.25
$ws.Range("K61").Value = 1498.25
$ws.Range("M61").Value = -1286.25

# Row 136 (G136 = 43999)
$ws.Range("H136").Value = 1598.2
$ws.Range("I136").Value = 1498.25
$ws.Range("K136").Value = 4494.75
$ws.Range("M136").Value = -1944.75

$ws = $wb.Worksheets.Item("BSM")
# Row 20 (G20 = 14149)
$ws.Range("H20").Value = 2995.3
$ws.Range("I20").Value = 2995.3
$ws.Range("K20").Value = 2995.3
$ws.Range("M20").Value = -2748.3

# Row 86 (G86 = 12526)
$ws.Range("H86").Value = 2216.5
$ws.Range("I86").Value = 2439.8
$ws.Range("J86").Value = 1100
$ws.Range("K86").Value = 2439.8
$ws.Range("L86").Value = 1100
$ws.Range("M86").Value = -1316.8
$ws.Range("N86").Value = -3346

# Row 89 (G89 = 12526)
$ws.Range("H89").Value = 2216.5
$ws.Range("I89").Value = 2439.8
$ws.Range("J89").Value = 1100
$ws.Range("K89").Value = 12199
$ws.Range("L89").Value = 5500
$ws.Range("M89").Value = -6583
$ws.Range("N89").Value = -16732

$ws = $wb.Worksheets.Item("CRP")
# Row 5 (G5 = 1893)
$ws.Range("H5").Value = 382.33334
$ws.Range("I5").Value = 382.33334
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 382.33334
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -270.33334
$ws.Range("N5").ClearContents()

# Row 22 (G22 = 5367)
$ws.Range("H22").Value = 70
$ws.Range("I22").Value = 70
$ws.Range("K22").Value = 70
$ws.Range("M22").Value = 280

# Row 132 (G132 = 44019)
$ws.Range("H132").Value = 1823.3334
$ws.Range("I132").Value = 1823.3334
$ws.Range("K132").Value = 5470.0002
$ws.Range("M132").Value = -2940.0002

$ws = $wb.Worksheets.Item("CUL")
# Row 5 (G5 = 43974)
$ws.Range("H5").Value = 2233.3333
$ws.Range("J5").Value = 1850
$ws.Range("L5").Value = 5550
$ws.Range("N5").Value = -5774

# Row 16 (G16 = 4641)
$ws.Range("H16").Value = 5000
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 5000
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 15000
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -15346

# Row 23 (G23 = 4858)
$ws.Range("H23").Value = 2159.7144
$ws.Range("I23").Value = 2023.6
$ws.Range("K23").Value = 6070.799999999999
$ws.Range("M23").Value = -5835.799999999999

# Row 50 (G50 = 4725)
$ws.Range("H50").Value = 661.2

# Row 53 (G53 = 4725)
$ws.Range("H53").Value = 661.2

# Row 92 (G92 = 19841)
$ws.Range("H92").Value = 5000
$ws.Range("I92").Value = 5000
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 15000
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = -13752
$ws.Range("N92").ClearContents()

# Row 97 (G97 = 19846)
$ws.Range("H97").Value = 3
$ws.Range("I97").Value = 3
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 9
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = 487
$ws.Range("N97").ClearContents()

# Row 107 (G107 = 27838)
$ws.Range("H107").Value = 507.14285
$ws.Range("I107").Value = 366.66666
$ws.Range("K107").Value = 1099.99998
$ws.Range("M107").Value = 820.0000199999999

# Row 121 (G121 = 27878)
$ws.Range("H121").Value = 1369.375
$ws.Range("I121").Value = 524.4
$ws.Range("J121").Value = 2777.6667
$ws.Range("K121").Value = 1573.2
$ws.Range("L121").Value = 8333.000100000001
$ws.Range("M121").Value = -263.1999999999998
$ws.Range("N121").Value = -10953.0001

# Row 135 (G135 = 43974)
$ws.Range("H135").Value = 2233.3333
$ws.Range("J135").Value = 1850
$ws.Range("L135").Value = 16650
$ws.Range("N135").Value = -21720

$ws = $wb.Worksheets.Item("GSM")
# Row 11 (G11 = 4422)
$ws.Range("H11").Value = 26183454
$ws.Range("I11").Value = 23500000
$ws.Range("J11").Value = 33339334
$ws.Range("K11").Value = 23500000
$ws.Range("L11").Value = 33339334
$ws.Range("M11").Value = -23499861
$ws.Range("N11").Value = -33339612

# Row 12 (G12 = 4093)
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()

# Row 97 (G97 = 19940)
$ws.Range("H97").Value = 1939.3334
$ws.Range("J97").Value = 2266
$ws.Range("L97").Value = 2266
$ws.Range("N97").Value = -3258

# Row 122 (G122 = 36182)
$ws.Range("H122").Value = 10355.947
$ws.Range("I122").Value = 6386.7144
$ws.Range("K122").Value = 19160.1432
$ws.Range("M122").Value = -16710.1432

$ws = $wb.Worksheets.Item("LTW")
# Row 32 (G32 = 2250)
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()

# Row 82 (G82 = 12565)
$ws.Range("H82").Value = 1999.3334
$ws.Range("I82").Value = 1999.3334
$ws.Range("K82").Value = 1999.3334
$ws.Range("M82").Value = -1638.3334

# Row 85 (G85 = 12565)
$ws.Range("H85").Value = 1999.3334
$ws.Range("I85").Value = 1999.3334
$ws.Range("K85").Value = 1999.3334
$ws.Range("M85").Value = -751.3334

$ws = $wb.Worksheets.Item("WVR")
# Row 132 (G132 = 44029)
$ws.Range("H132").Value = 5776.96
$ws.Range("I132").Value = 3045.625
$ws.Range("K132").Value = 9136.875
$ws.Range("M132").Value = -6606.875
